$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Build style palette (temporary helper cells) ---
$ws.Range("A1").Copy()
$ws.Range("Z101").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("Z102").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("Z103").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("Z104").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("Z105").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("Z106").PasteSpecial(-4122)
$ws.Range("B25").Copy()
$ws.Range("Z107").PasteSpecial(-4122)
$ws.Range("Z108").Interior.Color = 15453831
$ws.Range("Z108").Borders.LineStyle = 1
$ws.Range("Z109").Interior.Color = 42495
$ws.Range("Z109").Borders.LineStyle = 1

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 7.17
$ws.Columns.Item(5).ColumnWidth = 11.17
$ws.Columns.Item(6).ColumnWidth = 6.17

# --- Set cell values ---
$ws.Range("A1").Value = "Ticker"
$ws.Range("B1").Value = "ROC"
$ws.Range("C1").Value = "Mensual"
$ws.Range("D1").Value = "Semanal"
$ws.Range("E1").Value = "Trimestral"
$ws.Range("F1").Value = "Señal"
$ws.Range("A2").Value = "IBIT"
$ws.Range("B2").Value = 36.02
$ws.Range("A3").Value = "GLD"
$ws.Range("B3").Value = 23.45
$ws.Range("A4").Value = "FXI"
$ws.Range("B4").Value = 7.88
$ws.Range("A5").Value = "XLP"
$ws.Range("B5").Value = 2.87
$ws.Range("A6").Value = "XLF"
$ws.Range("B6").Value = 1.77
$ws.Range("A7").Value = "XLU"
$ws.Range("B7").Value = -1.6
$ws.Range("A8").Value = "SLV"
$ws.Range("B8").Value = -2.2
$ws.Range("A9").Value = "UUP"
$ws.Range("B9").Value = -3.13
$ws.Range("A10").Value = "TLT"
$ws.Range("B10").Value = -3.64
$ws.Range("A11").Value = "MTUM"
$ws.Range("B11").Value = -4.06
$ws.Range("A12").Value = "USO"
$ws.Range("B12").Value = -5.4
$ws.Range("A13").Value = "XLY"
$ws.Range("B13").Value = -6.05
$ws.Range("A14").Value = "XLRE"
$ws.Range("B14").Value = -6.5
$ws.Range("A15").Value = "XLI"
$ws.Range("B15").Value = -6.85
$ws.Range("A16").Value = "GMF"
$ws.Range("B16").Value = -7.11
$ws.Range("A17").Value = "XLE"
$ws.Range("B17").Value = -7.21
$ws.Range("A18").Value = "RSP"
$ws.Range("B18").Value = -7.36
$ws.Range("A19").Value = "SPY"
$ws.Range("B19").Value = -8.07
$ws.Range("A20").Value = "SPYV"
$ws.Range("B20").Value = -8.23
$ws.Range("A21").Value = "XLV"
$ws.Range("B21").Value = -8.43
$ws.Range("A22").Value = "SPYG"
$ws.Range("B22").Value = -8.77
$ws.Range("A23").Value = "QQQ"
$ws.Range("B23").Value = -9.43
$ws.Range("A24").Value = "XLB"
$ws.Range("B24").Value = -12.35
$ws.Range("A25").Value = "MOAT"
$ws.Range("B25").Value = -12.8
$ws.Range("A26").Value = "IWN"
$ws.Range("B26").Value = -13.37
$ws.Range("A27").Value = "IWO"
$ws.Range("B27").Value = -14.09
$ws.Range("A28").Value = "XLK"
$ws.Range("B28").Value = -15.77

# --- Apply styles (paste formats from palette) ---
$ws.Range("Z101").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("Z101").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("Z101").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("Z101").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("Z101").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("Z101").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("Z102").Copy()
$ws.Range("A2:A28").PasteSpecial(-4122)
$ws.Range("Z102").Copy()
$ws.Range("F2:F4").PasteSpecial(-4122)
$ws.Range("Z102").Copy()
$ws.Range("F8:F9").PasteSpecial(-4122)
$ws.Range("Z102").Copy()
$ws.Range("F11:F15").PasteSpecial(-4122)
$ws.Range("Z102").Copy()
$ws.Range("F18:F20").PasteSpecial(-4122)
$ws.Range("Z102").Copy()
$ws.Range("F22:F28").PasteSpecial(-4122)
$ws.Range("Z103").Copy()
$ws.Range("B2:B6").PasteSpecial(-4122)
$ws.Range("Z104").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("Z104").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("Z104").Copy()
$ws.Range("C10:C11").PasteSpecial(-4122)
$ws.Range("Z104").Copy()
$ws.Range("C13:C14").PasteSpecial(-4122)
$ws.Range("Z104").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("Z105").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("Z105").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("Z105").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("Z105").Copy()
$ws.Range("C17:C28").PasteSpecial(-4122)
$ws.Range("Z105").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("Z105").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("Z105").Copy()
$ws.Range("D6:D28").PasteSpecial(-4122)
$ws.Range("Z105").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("Z105").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("Z105").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("Z106").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("Z106").Copy()
$ws.Range("C5:C8").PasteSpecial(-4122)
$ws.Range("Z106").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("Z106").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("Z106").Copy()
$ws.Range("E2:E3").PasteSpecial(-4122)
$ws.Range("Z106").Copy()
$ws.Range("E5:E9").PasteSpecial(-4122)
$ws.Range("Z106").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("Z106").Copy()
$ws.Range("E13:E28").PasteSpecial(-4122)
$ws.Range("Z107").Copy()
$ws.Range("B7:B28").PasteSpecial(-4122)
$ws.Range("Z108").Copy()
$ws.Range("F6:F7").PasteSpecial(-4122)
$ws.Range("Z108").Copy()
$ws.Range("F10").PasteSpecial(-4122)
$ws.Range("Z108").Copy()
$ws.Range("F16:F17").PasteSpecial(-4122)
$ws.Range("Z108").Copy()
$ws.Range("F21").PasteSpecial(-4122)
$ws.Range("Z109").Copy()
$ws.Range("F5").PasteSpecial(-4122)

# --- Clean up helper/palette cells ---
$ws.Range("Z101:Z109").Clear()

$excel.CutCopyMode = $false